# "Added a point light to the second scene"
#
# Rubric sheet: rows 14-25 are the "KEY FEATURES" section for the second
# level/scene. Row 15 ("Support for adding Point and/or Spotlight sources
# via Blender") is being marked Complete, and an explanatory note is being
# left in column E for that row. The note that used to sit at E21/E30
# (a placeholder "d") is removed from those rows - only the one on row 15
# remains, now with real text describing the point light that was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Mark "Support for adding Point and/or Spotlight sources via Blender" complete
$ws.Range("D15").Value = "Complete"

# Leave a note about the change on the same row
$ws.Range("E15").Value = "Note: Level 2 has a yellow light in the center"

# Remove the stray placeholder notes that used to live on these rows
$ws.Range("E21").ClearContents()
$ws.Range("E30").ClearContents()

# Update the on-screen selection/scroll position to match where the edit was made
$ws.Range("E15").Select()
